$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new activity entry, no date
$ws.Range("B18").Value = "Contemplation on how to approach what I have left"
$ws.Range("C18").Value = 0.5

# Row 19: new activity entry with date (copy the existing date-formatted style first)
$ws.Range("A14").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = (Get-Date -Year 2018 -Month 12 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B19").Value = "Started to rework questions 2 and 3 for Data Analysis. Looked into how to deal with outliers + things about linear regression."
$ws.Range("C19").Value = 1

# Update selection to match target state
$ws.Range("H20").Select()
